$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-26 Thursday" "2024-12-27 Friday"
Replace-Text "327×2=" "728×9="
Replace-Text "470×9=" "806×7="
Replace-Text "778×6=" "295×9="
Replace-Text "440×6=" "228×8="
Replace-Text "693×5=" "734×2="
Replace-Text "147×8=" "407×6="
Replace-Text "493×6=" "149×6="
Replace-Text "836×4=" "895×8="
Replace-Text "110×2=" "469×3="
Replace-Text "693×2=" "812×6="
Replace-Text "910×8=" "467×5="
Replace-Text "119×8=" "675×8="
Replace-Text "681×2=" "712×7="
Replace-Text "717×4=" "621×6="
Replace-Text "725×5=" "840×3="
Replace-Text "794×9=" "579×4="
Replace-Text "924×2=" "688×5="
Replace-Text "213×3=" "968×7="
Replace-Text "812×8=" "119×2="
Replace-Text "948×4=" "238×7="
Replace-Text "837×4=" "120×6="
Replace-Text "171×9=" "944×2="
Replace-Text "404×4=" "651×5="
Replace-Text "207×8=" "518×5="
Replace-Text "972×7=" "172×2="

Write-Output "Done"
